# Update column G ("K" = strikeouts) with newly regenerated values.
# The rest of the sheet (rows, other columns) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 2
    4  = 3
    5  = 6
    6  = 7
    7  = 5
    8  = 9
    9  = 4
    10 = 6
    11 = 6
    12 = 6
    13 = 7
    14 = 5
    15 = 2
    16 = 4
    17 = 3
    18 = 8
    19 = 6
    20 = 3
    21 = 7
    22 = 12
    23 = 11
    24 = 7
    25 = 7
    26 = 6
    27 = 8
    28 = 4
    29 = 5
    30 = 4
    31 = 9
    32 = 8
    33 = 3
    34 = 2
    35 = 9
    36 = 4
    37 = 5
    38 = 4
    39 = 2
    40 = 3
    41 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
